{"js": "// Locate the two relevant paragraphs (\"Previous\" field name and the\n// \"  emp.var.rate\" field name) in the data-schema bullet list, then:\n//   1. Move the `_GoBack` bookmark from the start of the \"emp.var.rate\"\n//      paragraph to the end of the \"Previous\" paragraph.\n//   2. Strip the stray leading double-space from the \"emp.var.rate\" text.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet previousPara = null;\nlet empVarPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (text === \"Previous\") {\n    previousPara = para;\n  } else if (text.trim() === \"emp.var.rate\" && text !== \"emp.var.rate\") {\n    empVarPara = para;\n  }\n}\n\nif (!previousPara || !empVarPara) {\n  throw new Error(\"Could not locate the 'Previous' / 'emp.var.rate' paragraphs\");\n}\n\n// Remove the old bookmark wherever it currently lives.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Re-create it collapsed at the end of the \"Previous\" paragraph.\npreviousPara.getRange(\"End\").insertBookmark(\"_GoBack\");\n\n// Fix the \"emp.var.rate\" run text (drop the leading two spaces).\nempVarPara.getRange(\"Whole\").insertText(\"emp.var.rate\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Move the `_GoBack` bookmark from the start of the \"emp.var.rate\" bullet\n# to the end of the \"Previous\" bullet, and strip the stray leading\n# double-space from the \"emp.var.rate\" text.\n$doc = $word.ActiveDocument\n\n# --- Locate the two relevant paragraphs by their exact text -----------------\n$previousPara = $null\nforeach ($p in $doc.Paragraphs) {\n    if ($p.Range.Text -eq \"Previous`r\") {\n        $previousPara = $p\n        break\n    }\n}\nif ($previousPara -eq $null) {\n    throw \"Could not find the 'Previous' paragraph\"\n}\n\n# --- Remove the existing `_GoBack` bookmark (wherever it currently is) ------\nif ($doc.Bookmarks.Exists(\"_GoBack\")) {\n    $doc.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- Re-create it collapsed right after \"revious\", before the paragraph mark.\n# A bookmark built directly from a collapsed (zero-length) Range positioned\n# immediately before a paragraph mark can't be targeted reliably, so insert a\n# throwaway one-character marker there, wrap the bookmark around it, then\n# delete the marker \u2014 leaving the bookmark collapsed in the right spot.\n$endOfText = $previousPara.Range.End - 1\n$marker = $doc.Range($endOfText, $endOfText)\n$marker.InsertAfter(\"~\")\n$markerRange = $doc.Range($endOfText, $endOfText + 1)\n$doc.Bookmarks.Add(\"_GoBack\", $markerRange)\n$markerRange2 = $doc.Range($endOfText, $endOfText + 1)\n$markerRange2.Text = \"\"\n\n# --- Fix the \"emp.var.rate\" text: drop the leading two spaces ---------------\n$find = $doc.Content.Find\n$find.ClearFormatting()\n$find.Text = \"  emp.var.rate\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"emp.var.rate\"\n$find.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2)\n"}
